{"js": "// Locate the \"Remote\" list item (the last bullet of the \"\u0421\u043f\u0438\u0441\u044a\u043a \u0441\n// \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u0438\" / materials list: Arduino uno, x3 Servo motors, Ir receiver,\n// Wires, Power supply (9V battery), Remote) and insert a new list\n// paragraph right after it containing the Thingiverse project link,\n// inheriting the same \"List Paragraph\" formatting.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Remote\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  const newPara = target.insertParagraph(\n    \"https://www.thingiverse.com/thing:616239\",\n    Word.InsertLocation.after\n  );\n  await context.sync();\n}\n", "ps1": "# Locate the \"Remote\" list item (the last bullet of the \"\u0421\u043f\u0438\u0441\u044a\u043a \u0441\n# \u043c\u0430\u0442\u0435\u0440\u0438\u0430\u043b\u0438\" / materials list: Arduino uno, x3 Servo motors, Ir receiver,\n# Wires, Power supply (9V battery), Remote) and insert a new list\n# paragraph right after it containing the Thingiverse project link,\n# inheriting the same \"List Paragraph\" formatting.\n$d = $word.ActiveDocument\n\n$findRange = $d.Content\n$findRange.Find.MatchWholeWord = $true\n$found = $findRange.Find.Execute(\"Remote\")\n\nif ($found) {\n  $para = $findRange.Paragraphs(1)\n  $para.Range.InsertParagraphAfter()\n  $newPara = $para.Next()\n  $newPara.Range.InsertAfter(\"https://www.thingiverse.com/thing:616239\")\n}\n"}
